# Generate Report for Handback
# - Marks the two localization files as "Handed back: in sync with en-US"
#   instead of "Ready for handoff" on every sheet that tracks status.
# - Populates the newly-tracked "Latest Target File" / "Latest Handback File"
#   columns (F/G) on the zh-cn and de-de sheets, each linked back to the
#   source .md file and the locale-specific .xlf file (same targets already
#   used by the existing Latest Handoff File columns).
# - Stamps the "Latest Handback DateTime" (column H) now that a handback
#   actually happened for de-de, and that zh-cn's existing placeholder date
#   effectively becomes a real handback date.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# zh-cn: Latest Target File (F) / Latest Handback File (G) for both rows,
# hyperlinked the same way the existing Latest Handoff File (A) / Latest
# Handoff File (D) columns are.
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f3eb23229011dc97fb19a19d8024269a6ee95da7/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md",
    [Type]::Missing,
    [Type]::Missing,
    "d2f80547-3b07-445f-ae15-9c500b9db91d.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad925be95dc0b9e0e12e76e6756ebc9e085395cf/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f3eb23229011dc97fb19a19d8024269a6ee95da7/e2e/d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md",
    [Type]::Missing,
    [Type]::Missing,
    "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad925be95dc0b9e0e12e76e6756ebc9e085395cf/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.zh-cn.xlf"
) | Out-Null

# zh-cn had no real handback yet, but the placeholder date now reflects the
# actual (first) handback timestamp.
$wsZhCn.Range("H2").Value = "2016-03-23 09:38:16"
$wsZhCn.Range("H3").Value = "2016-03-23 09:38:16"

# de-de: same new Latest Target File (F) / Latest Handback File (G) columns.
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f3eb23229011dc97fb19a19d8024269a6ee95da7/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md",
    [Type]::Missing,
    [Type]::Missing,
    "d2f80547-3b07-445f-ae15-9c500b9db91d.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8960772060e52c3708254bb587723ac2fb69fd6b/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f3eb23229011dc97fb19a19d8024269a6ee95da7/e2e/d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md",
    [Type]::Missing,
    [Type]::Missing,
    "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8960772060e52c3708254bb587723ac2fb69fd6b/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "d4bd5148-b5f2-4a8e-aef0-31d565bbc1b7.f46e10b3ca3392df74d3aebdbec264278f4abd14.de-de.xlf"
) | Out-Null

# de-de actually got handed back - stamp the real handback datetime.
$wsDeDe.Range("H2").Value = "2016-03-23 09:38:29"
$wsDeDe.Range("H3").Value = "2016-03-23 09:38:29"
